# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# "Bad Drivers" table (rows 3-5 data, row 6 totals)
# ---------------------------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 366
$ws.Range("D3").Value = 96.2

$ws.Range("A4").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.110.0.5"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 83
$ws.Range("D4").Value = 97.3

$ws.Range("A5").Value = "Intel(R) Dual Band Wireless-AC 7260 - 17.15.0.5"
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 915
$ws.Range("D5").Value = 98.8

$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 1364

# ---------------------------------------------------------------
# "Good Drivers" table
# Rows 14-17 keep a blank "Driver Vintage" (now an empty numeric
# cell instead of a date string), rows 18-21 get new vintages, and
# three brand-new rows (22-24) are appended.
# ---------------------------------------------------------------

$ws.Range("A14").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.40.2.1"
$ws.Range("B14").Value = 27946
$ws.Range("D14").Value = 100
$ws.Range("E14").ClearContents()

$ws.Range("A15").Value = "Intel(R) Dual Band Wireless-AC 7260 - 18.33.15.1"
$ws.Range("B15").Value = 83189
$ws.Range("D15").Value = 100
$ws.Range("E15").ClearContents()

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B16").Value = 56018
$ws.Range("D16").Value = 100
$ws.Range("E16").ClearContents()

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B17").Value = 34244
$ws.Range("D17").Value = 100
$ws.Range("E17").ClearContents()

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B18").Value = 442178
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2024-11-10"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.60.1.2"
$ws.Range("B19").Value = 47426
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2024-06-02"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B20").Value = 77849
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2021-08-18"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B21").Value = 59673
$ws.Range("D21").Value = 100
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2020-08-05"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B22").Value = 113652
$ws.Range("D22").Value = 100
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2019-12-14"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.70.3.1"
$ws.Range("B23").Value = 20076
$ws.Range("D23").Value = 99.90000000000001
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2019-04-28"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.90.0.2"
$ws.Range("B24").Value = 53308
$ws.Range("D24").Value = 100
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2018-07-03"
